$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Right" (B) column values for the Marking and Total rows,
# and the corrected/total marks text for the Total row.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 75
$ws.Range("E12").Value = "75/140"
